# Populate the "Credit Card No" column (H) for each team row (2-29).
# Row 24 already had a credit card number (MLSC274581924053); per the
# target data it is replaced with a new number, and the rest of the
# rows (previously blank in column H) get newly assigned numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "MLSC274581924053"
$ws.Range("H3").Value = "MLSC273411206789"
$ws.Range("H4").Value = "MLSC278956012348"
$ws.Range("H5").Value = "MLSC271900439281"
$ws.Range("H6").Value = "MLSC276753908823"
$ws.Range("H7").Value = "MLSC278021677349"
$ws.Range("H8").Value = "MLSC279188325690"
$ws.Range("H9").Value = "MLSC274012093948"
$ws.Range("H10").Value = "MLSC273665718204"
$ws.Range("H11").Value = "MLSC279937456132"
$ws.Range("H12").Value = "MLSC275302947685"
$ws.Range("H13").Value = "MLSC272490411236"
$ws.Range("H14").Value = "MLSC278386074821"
$ws.Range("H15").Value = "MLSC271219486573"
$ws.Range("H16").Value = "MLSC275630089147"
$ws.Range("H17").Value = "MLSC277953712340"
$ws.Range("H18").Value = "MLSC273519849023"
$ws.Range("H19").Value = "MLSC272764021980"
$ws.Range("H20").Value = "MLSC278241857304"
$ws.Range("H21").Value = "MLSC279401358492"
$ws.Range("H22").Value = "MLSC275146789013"
$ws.Range("H23").Value = "MLSC273805276149"
$ws.Range("H24").Value = "MLSC271729503826"
$ws.Range("H25").Value = "MLSC276089314578"
$ws.Range("H26").Value = "MLSC274920348612"
$ws.Range("H27").Value = "MLSC273276041398"
$ws.Range("H28").Value = "MLSC275490028347"
$ws.Range("H29").Value = "MLSC277150283904"
